$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计" (i.e. right
#    before the current "2022-Q3" sheet), and populate it with the
#    fund holding breakdown for the new quarter.
# ------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

$headerLine = "基金代码	基金名称	基金规模	股票总仓位	仓位占比	持有市值(亿元)	仓位排名"
$headerCols = $headerLine -split "`t"
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $cell = $newSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headerCols[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$dataBlock = @"
004666	长城久嘉创新成长灵活配置混合A	24.38	92.75	5.07	1.2361	4
070002	嘉实增长混合	29.11	74.32	3.28	0.9548	8
010052	长城久嘉创新成长灵活配置混合C	15.65	92.75	5.07	0.7935	4
000654	华商新锐产业灵活配置混合	13.13	86.13	3.41	0.4477	4
002669	华商万众创新灵活配置混合A	11.81	89.02	3.41	0.4027	4
006803	嘉实互通精选股票	5.63	89.43	5.24	0.2950	3
010296	万家互联互通中国优势量化策略混合A	4.37	94.52	5.85	0.2556	4
970010	方正证券金立方一年持有期混合C	10.16	73.92	2.35	0.2388	8
012568	天弘高端制造混合A	6.15	90.96	3.60	0.2214	8
004423	华商研究精选灵活配置混合A	5.19	86.50	3.41	0.1770	4
003593	国泰景气行业灵活配置混合	3.51	92.21	4.48	0.1572	7
012491	华商核心引力混合A	3.88	89.14	3.41	0.1323	4
001760	嘉实创新成长灵活配置混合	1.39	82.45	8.48	0.1179	1
008961	华商科技创新混合	2.56	93.13	3.41	0.0873	3
016069	华商研究精选灵活配置混合C	1.59	86.50	3.41	0.0542	4
001758	嘉实研究增强灵活配置混合	0.96	93.60	5.01	0.0481	5
160722	嘉实惠泽灵活配置混合（LOF）	0.87	92.18	5.50	0.0478	2
000522	华润元大信息传媒科技混合	1.38	62.01	2.94	0.0406	10
012569	天弘高端制造混合C	1.08	90.96	3.60	0.0389	8
011369	华商均衡成长混合A	1.13	88.83	3.41	0.0385	4
010797	长城优选回报六个月持有期混合A	2.72	31.00	1.33	0.0362	2
016305	农银专精特新混合A	1.37	91.93	2.07	0.0284	9
010297	万家互联互通中国优势量化策略混合C	0.47	94.52	5.85	0.0275	4
002703	长城久源灵活配置混合A	0.70	94.29	3.75	0.0262	7
011370	华商均衡成长混合C	0.64	88.83	3.41	0.0218	4
002292	诺安益鑫灵活配置混合A	0.37	69.58	4.17	0.0154	7
012492	华商核心引力混合C	0.43	89.14	3.41	0.0147	4
350001	天治财富增长混合	0.53	69.74	2.58	0.0137	7
003238	新华外延增长主题灵活配置混合	0.51	85.17	2.44	0.0124	6
005088	嘉实新添辉定期开放灵活配置混合A	0.48	90.69	2.47	0.0119	3
006522	财通新兴蓝筹混合A	0.23	80.41	4.65	0.0107	3
004931	华润元大价值优选混合C	0.16	74.11	4.25	0.0068	8
009719	招商增浩一年定期开放混合C	0.71	21.56	0.94	0.0067	10
010798	长城优选回报六个月持有期混合C	0.41	31.00	1.33	0.0055	2
004930	华润元大价值优选混合A	0.13	74.11	4.25	0.0055	8
006523	财通新兴蓝筹混合C	0.11	80.41	4.65	0.0051	3
016306	农银专精特新混合C	0.21	91.93	2.07	0.0043	9
009718	招商增浩一年定期开放混合A	0.38	21.56	0.94	0.0036	10
014381	长城久源灵活配置混合C	0.06	94.29	3.75	0.0022	7
007875	国融融兴灵活配置混合A	0.03	73.25	4.71	0.0014	4
016051	华商万众创新灵活配置混合C	0.03	89.02	3.41	0.0010	4
970009	方正证券金立方一年持有期混合A	0.04	73.92	2.35	0.0009	8
014550	诺安益鑫灵活配置混合C	0.02	69.58	4.17	0.0008	7
007876	国融融兴灵活配置混合C	0.01	73.25	4.71	0.0005	4
005089	嘉实新添辉定期开放灵活配置混合C	0.01	90.69	2.47	0.0002	3
"@

$lines = $dataBlock -split "`n"
$rowIndex = 2
$idx = 0
foreach ($line in $lines) {
    $line = $line.Trim("`r")
    if ($line.Length -eq 0) { continue }
    $fields = $line -split "`t"
    $code = $fields[0]
    $name = $fields[1]
    $scale = $fields[2]
    $stockPos = $fields[3]
    $posRatio = $fields[4]
    $mktValue = $fields[5]
    $rank = $fields[6]

    $aCell = $newSheet.Cells.Item($rowIndex, 1)
    $aCell.Value = $idx
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $newSheet.Cells.Item($rowIndex, 2).Value = "'" + $code
    $newSheet.Cells.Item($rowIndex, 3).Value = "'" + $name
    $newSheet.Cells.Item($rowIndex, 4).Value = "'" + $scale
    $newSheet.Cells.Item($rowIndex, 5).Value = "'" + $stockPos
    $newSheet.Cells.Item($rowIndex, 6).Value = "'" + $posRatio
    $newSheet.Cells.Item($rowIndex, 7).Value = "'" + $mktValue
    $newSheet.Cells.Item($rowIndex, 8).Value = [double]$rank

    $rowIndex = $rowIndex + 1
    $idx = $idx + 1
}

Write-Host "2022-Q4 sheet populated with $($idx) rows"

# ------------------------------------------------------------------
# 2) Insert a new row 2 into "总计" with the 2022-Q4 summary, shifting
#    the existing quarters down by one row.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 45
$totalSheet.Cells.Item(2, 4).Value = 6.05

# Pick up the "s=2" look (thin border, centred) used by the rest of
# column A by copying the format from the cell directly below.
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

# The row-insert preserved the original (now stale) row-index values in
# column A for rows 3-8 -- renumber them sequentially (0-based) to match
# their new row position.
for ($r = 3; $r -le 8; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

Write-Host "总计 sheet row inserted"
